# Insert a new row at position 12 (shifts existing rows 12..121 down to 13..122)
# and populate it with the new Arveja Verde record described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 12.
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with the new record's data.
$ws.Cells.Item(12, 1).Value = 4
$ws.Cells.Item(12, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(12, 3).Value = "Los Lagos"
$ws.Cells.Item(12, 4).Value = 44750
$ws.Cells.Item(12, 5).Value = 10
$ws.Cells.Item(12, 6).Value = 100112022
$ws.Cells.Item(12, 7).Value = "Arveja Verde"
$ws.Cells.Item(12, 8).Value = "Perfection"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 70
$ws.Cells.Item(12, 11).Value = 40000
$ws.Cells.Item(12, 12).Value = 40000
$ws.Cells.Item(12, 13).Value = 40000
$ws.Cells.Item(12, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(12, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(12, 16).Value = 1600
$ws.Cells.Item(12, 17).Value = 25
$ws.Cells.Item(12, 18).Value = "Hortaliza"
